$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Worksheet")

$ws.Range("F1").Value = "% of Start"
$ws.Range("G1").Value = "% of Subpopulation"

$ws.Range("F1").Font.Bold = $true
$ws.Range("F1").NumberFormat = "0.00%"
$ws.Range("F1").WrapText = $true
$ws.Range("F1").Borders.LineStyle = 1

$ws.Range("G1").Font.Bold = $true
$ws.Range("G1").WrapText = $true
$ws.Range("G1").Borders.LineStyle = 1

$ws.Range("F3").Formula = "=C3/A$2"
$ws.Range("G3").Formula = "=C3/A$3"
$ws.Range("F4").Formula = "=C4/A$2"
$ws.Range("G4").Formula = "=C4/A$3"

$rng = $ws.Range("F3:G4")
$rng.NumberFormat = "0.00%"
$rng.WrapText = $true
$rng.Borders.LineStyle = 1
Write-Host "done"
